# Auto-generated script to update currentAveragePrice / LevePrice / LeveProfit
# columns (H..N) across all 8 item-category worksheets, reflecting a refreshed
# market-price data pull (chore: update Sheets via scheduled runner).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 27233.75
$ws.Range("J17").Value = 27233.75
$ws.Range("L17").Value = 81701.25
$ws.Range("N17").Value = -82037.25
$ws.Range("H18").Value = 12822779.0
$ws.Range("I18").Value = 12822779.0
$ws.Range("K18").Value = 12822779.0
$ws.Range("M18").Value = -12822495.0
$ws.Range("H19").Value = 1673.125
$ws.Range("I19").Value = 423.75
$ws.Range("J19").Value = 2922.5
$ws.Range("K19").Value = 423.75
$ws.Range("L19").Value = 2922.5
$ws.Range("M19").Value = -248.75
$ws.Range("N19").Value = -3272.5
$ws.Range("H28").Value = 662.4
$ws.Range("J28").Value = 2049.6667
$ws.Range("L28").Value = 2049.6667
$ws.Range("N28").Value = -3019.6667
$ws.Range("H32").Value = 3072.4473
$ws.Range("I32").Value = 2691.2856
$ws.Range("J32").Value = 3158.516
$ws.Range("K32").Value = 2691.2856
$ws.Range("L32").Value = 3158.516
$ws.Range("M32").Value = -2365.2856
$ws.Range("N32").Value = -3810.516
$ws.Range("H33").Value = 1372257.0
$ws.Range("I33").Value = 2315368.5
$ws.Range("J33").Value = 458.54544
$ws.Range("K33").Value = 2315368.5
$ws.Range("L33").Value = 458.54544
$ws.Range("M33").Value = -2315139.5
$ws.Range("N33").Value = -916.54544
$ws.Range("H53").Value = 6700.25
$ws.Range("I53").Value = 475.125
$ws.Range("J53").Value = 12925.375
$ws.Range("K53").Value = 475.125
$ws.Range("L53").Value = 12925.375
$ws.Range("M53").Value = 161.875
$ws.Range("N53").Value = -14199.375
$ws.Range("H62").Value = 4928.3076
$ws.Range("I62").Value = 3155.5715
$ws.Range("J62").Value = 6996.5
$ws.Range("K62").Value = 3155.5715
$ws.Range("L62").Value = 6996.5
$ws.Range("M62").Value = -2531.5715
$ws.Range("N62").Value = -8244.5
$ws.Range("H64").Value = 8203.889
$ws.Range("J64").Value = 8335.115
$ws.Range("L64").Value = 8335.115
$ws.Range("N64").Value = -8831.115
$ws.Range("H65").Value = 4928.3076
$ws.Range("I65").Value = 3155.5715
$ws.Range("J65").Value = 6996.5
$ws.Range("K65").Value = 15777.8575
$ws.Range("L65").Value = 34982.5
$ws.Range("M65").Value = -12657.8575
$ws.Range("N65").Value = -41222.5
$ws.Range("H67").Value = 8203.889
$ws.Range("J67").Value = 8335.115
$ws.Range("L67").Value = 8335.115
$ws.Range("N67").Value = -10051.115
$ws.Range("H76").Value = 5291.8887
$ws.Range("I76").Value = 4143.778
$ws.Range("J76").Value = 6440.0
$ws.Range("K76").Value = 4143.778
$ws.Range("L76").Value = 6440.0
$ws.Range("M76").Value = -3828.778
$ws.Range("N76").Value = -7070.0
$ws.Range("H79").Value = 5291.8887
$ws.Range("I79").Value = 4143.778
$ws.Range("J79").Value = 6440.0
$ws.Range("K79").Value = 4143.778
$ws.Range("L79").Value = 6440.0
$ws.Range("M79").Value = -3051.778
$ws.Range("N79").Value = -8624.0
$ws.Range("H86").Value = 5062.625
$ws.Range("I86").Value = 876.6667
$ws.Range("J86").Value = 7574.2
$ws.Range("K86").Value = 876.6667
$ws.Range("L86").Value = 7574.2
$ws.Range("M86").Value = 246.3333
$ws.Range("N86").Value = -9820.2
$ws.Range("H89").Value = 5062.625
$ws.Range("I89").Value = 876.6667
$ws.Range("J89").Value = 7574.2
$ws.Range("K89").Value = 4383.3335
$ws.Range("L89").Value = 37871.0
$ws.Range("M89").Value = 1232.6665
$ws.Range("N89").Value = -49103.0
$ws.Range("H92").Value = 950.62964
$ws.Range("I92").Value = 198.36363
$ws.Range("K92").Value = 198.36363
$ws.Range("M92").Value = 1049.63637
$ws.Range("H113").Value = 7095.2383
$ws.Range("I113").Value = 6711.3335
$ws.Range("K113").Value = 6711.3335
$ws.Range("M113").Value = -3457.3335
$ws.Range("H118").Value = 111111560.0
$ws.Range("I118").Value = 142857570.0
$ws.Range("K118").Value = 428572710.0
$ws.Range("M118").Value = -428571053.0
$ws.Range("H127").Value = 622.125
$ws.Range("I127").Value = 546.5
$ws.Range("J127").Value = 849.0
$ws.Range("K127").Value = 1639.5
$ws.Range("L127").Value = 2547.0
$ws.Range("M127").Value = 3320.5
$ws.Range("N127").Value = -12467.0
$ws.Range("H132").Value = 23259232.0
$ws.Range("I132").Value = 27781304.0
$ws.Range("J132").Value = 2855.0
$ws.Range("K132").Value = 83343912.0
$ws.Range("L132").Value = 8565.0
$ws.Range("M132").Value = -83341382.0
$ws.Range("N132").Value = -13625.0
$ws.Range("H135").Value = 1164.36
$ws.Range("I135").Value = 820.8095
$ws.Range("K135").Value = 7387.2855
$ws.Range("M135").Value = -4852.2855
$ws.Range("H137").Value = 54348.824
$ws.Range("I137").Value = 90379.0
$ws.Range("J137").Value = 2877.1428
$ws.Range("K137").Value = 271137.0
$ws.Range("L137").Value = 8631.4284
$ws.Range("M137").Value = -268587.0
$ws.Range("N137").Value = -13731.4284
$ws.Range("H138").Value = 3590.0303
$ws.Range("I138").Value = 2000.0
$ws.Range("J138").Value = 3692.6128
$ws.Range("K138").Value = 6000.0
$ws.Range("L138").Value = 11077.8384
$ws.Range("M138").Value = -860.0
$ws.Range("N138").Value = -21357.8384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6258026.0
$ws.Range("I45").Value = 11067534.0
$ws.Range("K45").Value = 11067534.0
$ws.Range("M45").Value = -11067157.0
$ws.Range("H74").Value = 36872.535
$ws.Range("I74").Value = 4941.72
$ws.Range("K74").Value = 4941.72
$ws.Range("M74").Value = -4067.72
$ws.Range("H77").Value = 36872.535
$ws.Range("I77").Value = 4941.72
$ws.Range("K77").Value = 24708.6
$ws.Range("M77").Value = -20340.6
$ws.Range("H88").Value = 1368.6923
$ws.Range("I88").Value = 1549.875
$ws.Range("K88").Value = 1549.875
$ws.Range("M88").Value = -1143.875
$ws.Range("H91").Value = 1368.6923
$ws.Range("I91").Value = 1549.875
$ws.Range("K91").Value = 1549.875
$ws.Range("M91").Value = -145.875
$ws.Range("H102").Value = 5955215.0
$ws.Range("I102").Value = 6946584.0
$ws.Range("J102").Value = 7000.0
$ws.Range("K102").Value = 6946584.0
$ws.Range("L102").Value = 7000.0
$ws.Range("M102").Value = -6944962.0
$ws.Range("N102").Value = -10244.0
$ws.Range("H122").Value = 871746.6
$ws.Range("I122").Value = 1748.0769
$ws.Range("J122").Value = 1899926.8
$ws.Range("K122").Value = 5244.2307
$ws.Range("L122").Value = 5699780.4
$ws.Range("M122").Value = -2794.2307
$ws.Range("N122").Value = -5704680.4
$ws.Range("H125").Value = 0.0
$ws.Range("J125").Value = 0.0
$ws.Range("L125").Value = 0.0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 3961.0
$ws.Range("I132").Value = 2935.75
$ws.Range("K132").Value = 8807.25
$ws.Range("M132").Value = -6277.25
$ws.Range("H140").Value = 61666.5
$ws.Range("J140").Value = 74999.5
$ws.Range("L140").Value = 74999.5
$ws.Range("N140").Value = -85359.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 1424.0
$ws.Range("I19").Value = 954.5
$ws.Range("J19").Value = 1893.5
$ws.Range("K19").Value = 954.5
$ws.Range("L19").Value = 1893.5
$ws.Range("M19").Value = -781.5
$ws.Range("N19").Value = -2239.5
$ws.Range("H86").Value = 5002721.5
$ws.Range("I86").Value = 7145605.0
$ws.Range("J86").Value = 2661.0
$ws.Range("K86").Value = 7145605.0
$ws.Range("L86").Value = 2661.0
$ws.Range("M86").Value = -7144482.0
$ws.Range("N86").Value = -4907.0
$ws.Range("H89").Value = 5002721.5
$ws.Range("I89").Value = 7145605.0
$ws.Range("J89").Value = 2661.0
$ws.Range("K89").Value = 35728025.0
$ws.Range("L89").Value = 13305.0
$ws.Range("M89").Value = -35722409.0
$ws.Range("N89").Value = -24537.0
$ws.Range("H94").Value = 3337104.5
$ws.Range("I94").Value = 3449418.2
$ws.Range("K94").Value = 3449418.2
$ws.Range("M94").Value = -3448967.2
$ws.Range("H105").Value = 3907201.8
$ws.Range("I105").Value = 3907201.8
$ws.Range("K105").Value = 3907201.8
$ws.Range("M105").Value = -3905454.8
$ws.Range("H107").Value = 8930149.0
$ws.Range("I107").Value = 8930149.0
$ws.Range("K107").Value = 8930149.0
$ws.Range("M107").Value = -8928229.0
$ws.Range("H134").Value = 6170.278
$ws.Range("I134").Value = 2008.25
$ws.Range("K134").Value = 6024.75
$ws.Range("M134").Value = -3489.75
$ws.Range("H138").Value = 64777.668
$ws.Range("J138").Value = 64777.668
$ws.Range("L138").Value = 64777.668
$ws.Range("N138").Value = -75057.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 120.388885
$ws.Range("I7").Value = 122.545456
$ws.Range("K7").Value = 122.545456
$ws.Range("M7").Value = -9.545456000000001
$ws.Range("H12").Value = 2002.5
$ws.Range("I12").Value = 2002.5
$ws.Range("K12").Value = 2002.5
$ws.Range("M12").Value = -1832.5
$ws.Range("H16").Value = 1845.9
$ws.Range("I16").Value = 1366.0
$ws.Range("J16").Value = 2165.8333
$ws.Range("K16").Value = 1366.0
$ws.Range("L16").Value = 2165.8333
$ws.Range("M16").Value = -1079.0
$ws.Range("N16").Value = -2739.8333
$ws.Range("H29").Value = 30333.0
$ws.Range("J29").Value = 30333.0
$ws.Range("L29").Value = 30333.0
$ws.Range("N29").Value = -30919.0
$ws.Range("H62").Value = 4099.5
$ws.Range("I62").Value = 3450.0
$ws.Range("J62").Value = 4749.0
$ws.Range("K62").Value = 3450.0
$ws.Range("L62").Value = 4749.0
$ws.Range("M62").Value = -2826.0
$ws.Range("N62").Value = -5997.0
$ws.Range("H65").Value = 4099.5
$ws.Range("I65").Value = 3450.0
$ws.Range("J65").Value = 4749.0
$ws.Range("K65").Value = 17250.0
$ws.Range("L65").Value = 23745.0
$ws.Range("M65").Value = -14130.0
$ws.Range("N65").Value = -29985.0
$ws.Range("H103").Value = 7000.25
$ws.Range("J103").Value = 20577.0
$ws.Range("L103").Value = 20577.0
$ws.Range("N103").Value = -22921.0
$ws.Range("H113").Value = 1845.9
$ws.Range("I113").Value = 1366.0
$ws.Range("J113").Value = 2165.8333
$ws.Range("K113").Value = 1366.0
$ws.Range("L113").Value = 2165.8333
$ws.Range("M113").Value = 804.0
$ws.Range("N113").Value = -6505.8333
$ws.Range("H132").Value = 92484.86
$ws.Range("I132").Value = 68555.6
$ws.Range("K132").Value = 205666.8
$ws.Range("M132").Value = -203136.8
$ws.Range("H134").Value = 3351.0908
$ws.Range("I134").Value = 2601.5789
$ws.Range("K134").Value = 7804.736699999999
$ws.Range("M134").Value = -5269.736699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 537.9583
$ws.Range("J23").Value = 612.7368
$ws.Range("L23").Value = 1838.2104
$ws.Range("N23").Value = -2308.2104
$ws.Range("H34").Value = 630.0
$ws.Range("I34").Value = 450.0
$ws.Range("J34").Value = 750.0
$ws.Range("K34").Value = 1350.0
$ws.Range("L34").Value = 2250.0
$ws.Range("M34").Value = -1266.0
$ws.Range("N34").Value = -2418.0
$ws.Range("H37").Value = 52737.5
$ws.Range("J37").Value = 52737.5
$ws.Range("L37").Value = 158212.5
$ws.Range("N37").Value = -158436.5
$ws.Range("H38").Value = 85.84615
$ws.Range("J38").Value = 183.25
$ws.Range("L38").Value = 549.75
$ws.Range("N38").Value = -1243.75
$ws.Range("H39").Value = 0.0
$ws.Range("I39").Value = 0.0
$ws.Range("J39").Value = 0.0
$ws.Range("K39").Value = 0.0
$ws.Range("L39").Value = 0.0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H57").Value = 3259.8
$ws.Range("I57").Value = 649.5
$ws.Range("J57").Value = 5000.0
$ws.Range("K57").Value = 1948.5
$ws.Range("L57").Value = 15000.0
$ws.Range("M57").Value = -1389.5
$ws.Range("N57").Value = -16118.0
$ws.Range("H58").Value = 1607.4286
$ws.Range("I58").Value = 752.5
$ws.Range("J58").Value = 1749.9166
$ws.Range("K58").Value = 2257.5
$ws.Range("L58").Value = 5249.7498
$ws.Range("M58").Value = -2129.5
$ws.Range("N58").Value = -5505.7498
$ws.Range("H98").Value = 1905.75
$ws.Range("J98").Value = 2049.4285
$ws.Range("L98").Value = 6148.2855
$ws.Range("N98").Value = -9144.2855
$ws.Range("H136").Value = 3373.6
$ws.Range("I136").Value = 3373.6
$ws.Range("K136").Value = 10120.8
$ws.Range("M136").Value = -5020.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 27527.0
$ws.Range("J32").Value = 27527.0
$ws.Range("L32").Value = 27527.0
$ws.Range("N32").Value = -28119.0
$ws.Range("H64").Value = 0.0
$ws.Range("J64").Value = 0.0
$ws.Range("L64").Value = 0.0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0.0
$ws.Range("J67").Value = 0.0
$ws.Range("L67").Value = 0.0
$ws.Range("N67").ClearContents()
$ws.Range("H80").Value = 37459110.0
$ws.Range("I80").Value = 47674330.0
$ws.Range("J80").Value = 3332.6667
$ws.Range("K80").Value = 47674330.0
$ws.Range("L80").Value = 3332.6667
$ws.Range("M80").Value = -47673332.0
$ws.Range("N80").Value = -5328.6667
$ws.Range("H83").Value = 37459110.0
$ws.Range("I83").Value = 47674330.0
$ws.Range("J83").Value = 3332.6667
$ws.Range("K83").Value = 238371650.0
$ws.Range("L83").Value = 16663.3335
$ws.Range("M83").Value = -238366658.0
$ws.Range("N83").Value = -26647.3335
$ws.Range("H99").Value = 16848.0
$ws.Range("I99").Value = 16113.8
$ws.Range("J99").Value = 20519.0
$ws.Range("K99").Value = 16113.8
$ws.Range("L99").Value = 20519.0
$ws.Range("M99").Value = -13867.8
$ws.Range("N99").Value = -25011.0
$ws.Range("H102").Value = 5624145.0
$ws.Range("I102").Value = 8548115.0
$ws.Range("K102").Value = 8548115.0
$ws.Range("M102").Value = -8546493.0
$ws.Range("H122").Value = 299852.38
$ws.Range("I122").Value = 389157.8
$ws.Range("K122").Value = 1167473.4
$ws.Range("M122").Value = -1165023.4
$ws.Range("H126").Value = 4956398.0
$ws.Range("I126").Value = 4547920.5
$ws.Range("J126").Value = 5211696.5
$ws.Range("K126").Value = 13643761.5
$ws.Range("L126").Value = 15635089.5
$ws.Range("M126").Value = -13641291.5
$ws.Range("N126").Value = -15640029.5
$ws.Range("H132").Value = 3777.3333
$ws.Range("I132").Value = 3726.318
$ws.Range("J132").Value = 4001.8
$ws.Range("K132").Value = 11178.954
$ws.Range("L132").Value = 12005.4
$ws.Range("M132").Value = -8648.954000000002
$ws.Range("N132").Value = -17065.4
$ws.Range("H134").Value = 39695.0
$ws.Range("J134").Value = 39695.0
$ws.Range("L134").Value = 119085.0
$ws.Range("N134").Value = -124155.0
$ws.Range("H141").Value = 50171.4
$ws.Range("J141").Value = 50171.4
$ws.Range("L141").Value = 50171.4
$ws.Range("N141").Value = -60531.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6800.875
$ws.Range("I7").Value = 4330.8184
$ws.Range("J7").Value = 8890.923
$ws.Range("K7").Value = 4330.8184
$ws.Range("L7").Value = 8890.923
$ws.Range("M7").Value = -4218.8184
$ws.Range("N7").Value = -9114.923
$ws.Range("H22").Value = 43717.715
$ws.Range("J22").Value = 2533.0
$ws.Range("L22").Value = 2533.0
$ws.Range("N22").Value = -3123.0
$ws.Range("H27").Value = 43717.715
$ws.Range("J27").Value = 2533.0
$ws.Range("L27").Value = 2533.0
$ws.Range("N27").Value = -2747.0
$ws.Range("H61").Value = 15874016.0
$ws.Range("H68").Value = 4333.3335
$ws.Range("J68").Value = 4750.0
$ws.Range("L68").Value = 4750.0
$ws.Range("N68").Value = -6248.0
$ws.Range("H71").Value = 4333.3335
$ws.Range("J71").Value = 4750.0
$ws.Range("L71").Value = 23750.0
$ws.Range("N71").Value = -31238.0
$ws.Range("H82").Value = 52779776.0
$ws.Range("I82").Value = 95962130.0
$ws.Range("J82").Value = 1342.0
$ws.Range("K82").Value = 95962130.0
$ws.Range("L82").Value = 1342.0
$ws.Range("M82").Value = -95961769.0
$ws.Range("N82").Value = -2064.0
$ws.Range("H85").Value = 52779776.0
$ws.Range("I85").Value = 95962130.0
$ws.Range("J85").Value = 1342.0
$ws.Range("K85").Value = 95962130.0
$ws.Range("L85").Value = 1342.0
$ws.Range("M85").Value = -95960882.0
$ws.Range("N85").Value = -3838.0
$ws.Range("H93").Value = 17555878.0
$ws.Range("I93").Value = 25642752.0
$ws.Range("K93").Value = 25642752.0
$ws.Range("M93").Value = -25641504.0
$ws.Range("H100").Value = 2310.926
$ws.Range("I100").Value = 1871.3572
$ws.Range("J100").Value = 2784.3076
$ws.Range("K100").Value = 1871.3572
$ws.Range("L100").Value = 2784.3076
$ws.Range("M100").Value = -1330.3572
$ws.Range("N100").Value = -3866.3076
$ws.Range("H108").Value = 0.0
$ws.Range("J108").Value = 0.0
$ws.Range("L108").Value = 0.0
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 0.0
$ws.Range("J109").Value = 0.0
$ws.Range("L109").Value = 0.0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 15874016.0
$ws.Range("H122").Value = 7161.6665
$ws.Range("I122").Value = 4555.5
$ws.Range("K122").Value = 13666.5
$ws.Range("M122").Value = -11216.5
$ws.Range("H126").Value = 6800.875
$ws.Range("I126").Value = 4330.8184
$ws.Range("J126").Value = 8890.923
$ws.Range("K126").Value = 12992.4552
$ws.Range("L126").Value = 26672.769
$ws.Range("M126").Value = -10522.4552
$ws.Range("N126").Value = -31612.769
$ws.Range("H136").Value = 29136.844
$ws.Range("I136").Value = 43008.57
$ws.Range("K136").Value = 129025.71
$ws.Range("M136").Value = -126475.71

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 34000.0
$ws.Range("J40").Value = 34000.0
$ws.Range("L40").Value = 34000.0
$ws.Range("N40").Value = -34298.0
$ws.Range("H93").Value = 60400.0
$ws.Range("J93").Value = 60400.0
$ws.Range("L93").Value = 60400.0
$ws.Range("N93").Value = -65392.0
$ws.Range("H95").Value = 0.0
$ws.Range("J95").Value = 0.0
$ws.Range("L95").Value = 0.0
$ws.Range("N95").ClearContents()
$ws.Range("H100").Value = 857.0714
$ws.Range("I100").Value = 1144.7142
$ws.Range("J100").Value = 569.4286
$ws.Range("K100").Value = 2289.4284
$ws.Range("L100").Value = 1138.8572
$ws.Range("M100").Value = -1748.4284
$ws.Range("N100").Value = -2220.8572
$ws.Range("H132").Value = 40442264.0
$ws.Range("I132").Value = 52639804.0
$ws.Range("J132").Value = 1816730.6
$ws.Range("K132").Value = 157919412.0
$ws.Range("L132").Value = 5450191.800000001
$ws.Range("M132").Value = -157916882.0
$ws.Range("N132").Value = -5455251.800000001
$ws.Range("H136").Value = 3358.1428
$ws.Range("I136").Value = 3511.5
$ws.Range("K136").Value = 10534.5
$ws.Range("M136").Value = -7984.5
$ws.Range("H140").Value = 0.0
$ws.Range("J140").Value = 0.0
$ws.Range("L140").Value = 0.0
$ws.Range("N140").ClearContents()

